# bioSample_hbrown_01.29.20.xlsx — "continuing major accuracy cleaning"
#
# The floodmedia column (H) was mis-populated with "NA" for every data row;
# it should read "None" instead (inductionDelay in column I correctly stays
# "NA"). Also nudge a couple of cosmetic view bits (row height, a couple of
# column widths, and the active-cell selection) to match the follow-up
# touch-up the author made while in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core data fix: floodmedia (column H) "NA" -> "None" for rows 2-27 ---
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 8).Value = "None"
}

# --- Cosmetic touch-ups matching the rest of the author's edit ---

# Data rows got a bit shorter (16pt -> 15pt); the header row (1) is untouched.
$ws.Range("A2:A27").RowHeight = 15

# Column I (inductionDelay) keeps ~the same width; column J (treatment) widens.
$ws.Columns.Item(9).ColumnWidth = 16.3
$ws.Columns.Item(10).ColumnWidth = 20

# Cursor ends up parked on I2.
$ws.Range("I2").Select()
